$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.761.26"
$ws.Range("E2").Value = "  -2.11%  "

$ws.Range("D3").Value = "3.341.95"
$ws.Range("E3").Value = "  -2.45%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.95"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.22"
$ws.Range("E6").Value = "  -3.55%  "

$ws.Range("E7").Value = "  -3.12%  "

$ws.Range("D8").Value = "3.334.54"
$ws.Range("E8").Value = "  -2.50%  "

$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("E10").Value = "  -1.80%  "

$ws.Range("E11").Value = "  +1.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.75"
$ws.Range("E12").Value = "  +0.71%  "

$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.89"
$ws.Range("E14").Value = "  -2.59%  "

$ws.Range("D15").Value = "3.877.05"
$ws.Range("E15").Value = "  -2.28%  "

$ws.Range("D16").Value = "3.363.49"
$ws.Range("E16").Value = "  -1.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "17.90"
$ws.Range("E17").Value = "  -1.72%  "

$ws.Range("E18").Value = "  -3.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.70"
$ws.Range("E19").Value = "  -0.60%  "

$ws.Range("D20").Value = "63.712.43"
$ws.Range("E20").Value = "  -2.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.974"
$ws.Range("E21").Value = "  -0.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "411.64"
$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.03"
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.34"
$ws.Range("E24").Value = "  +1.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.63"
$ws.Range("E25").Value = "  +12.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "82.96"
$ws.Range("E26").Value = "  -1.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.53"
$ws.Range("E27").Value = "  -2.15%  "

$ws.Range("E28").Value = "  -3.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.59"
$ws.Range("E29").Value = "  -2.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.06"
$ws.Range("E30").Value = "  -1.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.37"
$ws.Range("E31").Value = "  -1.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.31"
$ws.Range("E32").Value = "  -2.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "562.63"
$ws.Range("E33").Value = "  -7.80%  "

$ws.Range("E34").Value = "  -1.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.98"
$ws.Range("E35").Value = "  -0.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.147"
$ws.Range("E36").Value = "  +0.53%  "

$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.04"
$ws.Range("E38").Value = "  -5.70%  "

$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0737"
$ws.Range("E39").Value = "  -4.55%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.37"
$ws.Range("E40").Value = "  +1.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.366"
$ws.Range("E41").Value = "  -2.62%  "

$ws.Range("D42").Value = "3.140.80"
$ws.Range("E42").Value = "  -0.54%  "

$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("E44").Value = "  +0.69%  "

$ws.Range("E45").Value = "  +1.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0399"
$ws.Range("E46").Value = "  -2.01%  "

$ws.Range("E47").Value = "  -4.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.60"
$ws.Range("E48").Value = "  -4.29%  "

$ws.Range("E49").Value = "  -2.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.18"
$ws.Range("E50").Value = "  -4.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.06"
$ws.Range("E51").Value = "  -3.28%  "
